$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data values (commit: Updated cryptos list on Wed Apr 12 05:58:21 UTC 2023 with GitHub Actions)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.982.02"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.62"
$ws.Range("E3").Value = "  -2.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.11"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5081"
$ws.Range("E7").Value = "  -3.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3937"
$ws.Range("E8").Value = "  -2.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08171"
$ws.Range("E9").Value = "  -3.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.09"
$ws.Range("E10").Value = "  -1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.091"
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.73"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.91"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.257"
$ws.Range("E14").Value = "  -1.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.168"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.82"
$ws.Range("E17").Value = "  -4.50%  "
$ws.Range("E18").Value = "  -2.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06412"
$ws.Range("E19").Value = "  -4.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.959.57"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.802"
$ws.Range("E23").Value = "  -4.15%  "
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.150"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.087.20"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.83"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -9.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.23"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.060"
$ws.Range("E31").Value = "  -1.83%  "
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.892"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.728"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02421"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.209"
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06349"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.480"
$ws.Range("E40").Value = "  -5.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6297"
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.20"
$ws.Range("E42").Value = "  -3.36%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.199"
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.00"
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5890"
$ws.Range("E46").Value = "  -4.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.627"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.50"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.126"
$ws.Range("E51").Value = "  -3.66%  "
